$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2098.875
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H62").Value = 7358115.5
$ws.Range("I62").Value = 10420664
$ws.Range("K62").Value = 10420664
$ws.Range("M62").Value = -10420040

$ws.Range("H65").Value = 7358115.5
$ws.Range("I65").Value = 10420664
$ws.Range("K65").Value = 52103320
$ws.Range("M65").Value = -52100200

$ws.Range("H113").Value = 4130.3076
$ws.Range("I113").Value = 3533.111
$ws.Range("K113").Value = 3533.111
$ws.Range("M113").Value = -279.1109999999999

$ws.Range("H116").Value = 13514.083
$ws.Range("I116").Value = 4652.8
$ws.Range("J116").Value = 19843.572
$ws.Range("K116").Value = 4652.8
$ws.Range("L116").Value = 19843.572
$ws.Range("M116").Value = -1210.8
$ws.Range("N116").Value = -26727.572

$ws.Range("H125").Value = 9263477
$ws.Range("I125").Value = 2196.3333
$ws.Range("K125").Value = 19766.9997
$ws.Range("M125").Value = -17306.9997

$ws.Range("H127").Value = 2671.8
$ws.Range("I127").Value = 464.875
$ws.Range("J127").Value = 11499.5
$ws.Range("K127").Value = 1394.625
$ws.Range("L127").Value = 34498.5
$ws.Range("M127").Value = 3565.375
$ws.Range("N127").Value = -44418.5

$ws.Range("H132").Value = 2543.15
$ws.Range("I132").Value = 3148.4666
$ws.Range("J132").Value = 727.2
$ws.Range("K132").Value = 9445.399800000001
$ws.Range("L132").Value = 2181.6
$ws.Range("M132").Value = -6915.399800000001
$ws.Range("N132").Value = -7241.6

$ws.Range("H135").Value = 770657.4399999999
$ws.Range("I135").Value = 1001201.3
$ws.Range("K135").Value = 9010811.700000001
$ws.Range("M135").Value = -9008276.700000001

$ws.Range("H137").Value = 2708.36
$ws.Range("J137").Value = 2743.2222
$ws.Range("L137").Value = 8229.6666
$ws.Range("N137").Value = -13329.6666

$ws.Range("H139").Value = 60833.332
$ws.Range("J139").Value = 60833.332
$ws.Range("L139").Value = 60833.332
$ws.Range("N139").Value = -71113.33199999999

$ws.Range("H141").Value = 3549
$ws.Range("I141").Value = 3549
$ws.Range("K141").Value = 10647
$ws.Range("M141").Value = -5467

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4082.475
$ws.Range("I32").Value = 4082.475
$ws.Range("K32").Value = 4082.475
$ws.Range("M32").Value = -3795.475

$ws.Range("H48").Value = 239684
$ws.Range("J48").Value = 239684
$ws.Range("L48").Value = 239684
$ws.Range("N48").Value = -240452

$ws.Range("H61").Value = 2186.0454
$ws.Range("I61").Value = 1956.8096
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 1956.8096
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -1744.8096
$ws.Range("N61").Value = -7424

$ws.Range("H97").Value = 1172.6154
$ws.Range("I97").Value = 1085.1666
$ws.Range("K97").Value = 1085.1666
$ws.Range("M97").Value = -589.1666

$ws.Range("H136").Value = 2186.0454
$ws.Range("I136").Value = 1956.8096
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 5870.4288
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -3320.4288
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1948.6
$ws.Range("I11").Value = 2421.8
$ws.Range("J11").Value = 1475.4
$ws.Range("K11").Value = 2421.8
$ws.Range("L11").Value = 1475.4
$ws.Range("M11").Value = -2281.8
$ws.Range("N11").Value = -1755.4

$ws.Range("H42").Value = 164683.5
$ws.Range("J42").Value = 164683.5
$ws.Range("L42").Value = 164683.5
$ws.Range("N42").Value = -165339.5

$ws.Range("H43").Value = 222684
$ws.Range("J43").Value = 222684
$ws.Range("L43").Value = 222684
$ws.Range("N43").Value = -223046

$ws.Range("H48").Value = 239684
$ws.Range("J48").Value = 239684
$ws.Range("L48").Value = 239684
$ws.Range("N48").Value = -240514

$ws.Range("H59").Value = 87694.5
$ws.Range("J59").Value = 87694.5
$ws.Range("L59").Value = 87694.5
$ws.Range("N59").Value = -89388.5

$ws.Range("H94").Value = 150
$ws.Range("I94").Value = 150
$ws.Range("K94").Value = 150
$ws.Range("M94").Value = 301

$ws.Range("H139").Value = 78307.5
$ws.Range("J139").Value = 78307.5
$ws.Range("L139").Value = 78307.5
$ws.Range("N139").Value = -88587.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43882.84
$ws.Range("J31").Value = 54118.9
$ws.Range("L31").Value = 54118.9
$ws.Range("N31").Value = -54708.9

$ws.Range("H34").Value = 43882.84
$ws.Range("J34").Value = 54118.9
$ws.Range("L34").Value = 54118.9
$ws.Range("N34").Value = -54522.9

$ws.Range("H105").Value = 467
$ws.Range("I105").Value = 482.125
$ws.Range("K105").Value = 482.125
$ws.Range("M105").Value = 1264.875

$ws.Range("H134").Value = 214597.4
$ws.Range("I134").Value = 1897.289
$ws.Range("K134").Value = 5691.867
$ws.Range("M134").Value = -3156.867

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 7813018.5
$ws.Range("I11").Value = 14706165
$ws.Range("K11").Value = 44118495
$ws.Range("M11").Value = -44118355

$ws.Range("H23").Value = 239.66667
$ws.Range("I23").Value = 121
$ws.Range("K23").Value = 363
$ws.Range("M23").Value = -128

$ws.Range("H132").Value = 1378309.1
$ws.Range("I132").Value = 334823.34
$ws.Range("K132").Value = 3013410.06
$ws.Range("M132").Value = -3010880.06

$ws.Range("H140").Value = 2622.16
$ws.Range("I140").Value = 1402.7
$ws.Range("J140").Value = 7500
$ws.Range("K140").Value = 4208.1
$ws.Range("L140").Value = 22500
$ws.Range("M140").Value = 971.8999999999996
$ws.Range("N140").Value = -32860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6737.0435
$ws.Range("I70").Value = 7403.4375
$ws.Range("K70").Value = 7403.4375
$ws.Range("M70").Value = -7133.4375

$ws.Range("H73").Value = 6737.0435
$ws.Range("I73").Value = 7403.4375
$ws.Range("K73").Value = 7403.4375
$ws.Range("M73").Value = -6467.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 840897.75
$ws.Range("I7").Value = 10100.6
$ws.Range("J7").Value = 1434324.2
$ws.Range("K7").Value = 10100.6
$ws.Range("L7").Value = 1434324.2
$ws.Range("M7").Value = -9988.6
$ws.Range("N7").Value = -1434548.2

$ws.Range("H40").Value = 117052.555
$ws.Range("I40").Value = 171578.83
$ws.Range("K40").Value = 171578.83
$ws.Range("M40").Value = -171442.83

$ws.Range("H61").Value = 5389.727
$ws.Range("I61").Value = 5536.625
$ws.Range("K61").Value = 5536.625
$ws.Range("M61").Value = -5334.625

$ws.Range("H93").Value = 4638.8
$ws.Range("I93").Value = 4499
$ws.Range("K93").Value = 4499
$ws.Range("M93").Value = -3251

$ws.Range("H113").Value = 5389.727
$ws.Range("I113").Value = 5536.625
$ws.Range("K113").Value = 5536.625
$ws.Range("M113").Value = -3366.625

$ws.Range("H126").Value = 840897.75
$ws.Range("I126").Value = 10100.6
$ws.Range("J126").Value = 1434324.2
$ws.Range("K126").Value = 30301.8
$ws.Range("L126").Value = 4302972.6
$ws.Range("M126").Value = -27831.8
$ws.Range("N126").Value = -4307912.6

$ws.Range("H137").Value = 53750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8265.429
$ws.Range("I74").Value = 4930.6665
$ws.Range("J74").Value = 9174.909
$ws.Range("K74").Value = 4930.6665
$ws.Range("L74").Value = 9174.909
$ws.Range("M74").Value = -3994.6665
$ws.Range("N74").Value = -11046.909

$ws.Range("H77").Value = 8265.429
$ws.Range("I77").Value = 4930.6665
$ws.Range("J77").Value = 9174.909
$ws.Range("K77").Value = 14791.9995
$ws.Range("L77").Value = 27524.727
$ws.Range("M77").Value = -10111.9995
$ws.Range("N77").Value = -36884.727

$ws.Range("H126").Value = 1643.125
$ws.Range("I126").Value = 1643.125
$ws.Range("K126").Value = 4929.375
$ws.Range("M126").Value = -2459.375

$ws.Range("H136").Value = 8424760
$ws.Range("I136").Value = 10406172
$ws.Range("J136").Value = 251437.25
$ws.Range("K136").Value = 31218516
$ws.Range("L136").Value = 754311.75
$ws.Range("M136").Value = -31215966
$ws.Range("N136").Value = -759411.75

$ws.Range("H138").Value = 500043260
$ws.Range("J138").Value = 500043260
$ws.Range("L138").Value = 500043260
$ws.Range("N138").Value = -500053540
